$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (recalculated figures) ---
$ws.Range("G4").Value = 9117.377889493921
$ws.Range("G5").Value = 9117.377889493921
$ws.Range("G6").Value = 9117.377889493921
$ws.Range("G7").Value = 9117.377889493921
$ws.Range("G8").Value = 2.13015787156702
$ws.Range("G9").Value = 2.13015787156702
$ws.Range("G10").Value = 0.477676758834549
$ws.Range("G11").Value = 0.477676758834549
$ws.Range("F12").Value = 2.0448
$ws.Range("G12").Value = 4.99122545454545
$ws.Range("H12").Value = 26.1187
$ws.Range("I12").Value = 23.07283
$ws.Range("M12").Value = 10.34987
$ws.Range("N12").Value = 17.54201
$ws.Range("F13").Value = 2.0448
$ws.Range("G13").Value = 4.99122545454545
$ws.Range("H13").Value = 26.1187
$ws.Range("I13").Value = 23.07283
$ws.Range("M13").Value = 10.34987
$ws.Range("N13").Value = 17.54201
$ws.Range("G20").Value = 5232.08698040301
$ws.Range("G21").Value = 5232.08698040301
$ws.Range("G22").Value = 5232.08698040301
$ws.Range("G23").Value = 5232.08698040301
$ws.Range("G24").Value = 3.09641660704427
$ws.Range("G25").Value = 3.09641660704427
$ws.Range("G26").Value = 0.383764620864631
$ws.Range("G27").Value = 0.383764620864631
$ws.Range("G28").Value = 7.04017454545455
$ws.Range("H28").Value = 26.1187
$ws.Range("I28").Value = 23.2695
$ws.Range("M28").Value = 15.80875
$ws.Range("N28").Value = 22.46187
$ws.Range("G29").Value = 7.04017454545455
$ws.Range("H29").Value = 26.1187
$ws.Range("I29").Value = 23.2695
$ws.Range("M29").Value = 15.80875
$ws.Range("N29").Value = 22.46187
$ws.Range("G36").Value = 4830.94257003867
$ws.Range("G37").Value = 4830.94257003867
$ws.Range("G38").Value = 4830.94257003867
$ws.Range("G39").Value = 4830.94257003867
$ws.Range("G42").Value = 0.46512735331533
$ws.Range("G43").Value = 0.46512735331533
$ws.Range("G44").Value = 9.388478571428569
$ws.Range("H44").Value = 26.1187
$ws.Range("I44").Value = 23.26286
$ws.Range("L44").Value = 13.8189
$ws.Range("N44").Value = 22.39893
$ws.Range("G45").Value = 9.388478571428569
$ws.Range("H45").Value = 26.1187
$ws.Range("I45").Value = 23.26286
$ws.Range("L45").Value = 13.8189
$ws.Range("N45").Value = 22.39893
$ws.Range("G52").Value = 7238.14805810142
$ws.Range("G53").Value = 7238.14805810142
$ws.Range("G54").Value = 7238.14805810142
$ws.Range("G55").Value = 7238.14805810142
$ws.Range("G58").Value = 0.384581259397517
$ws.Range("G59").Value = 0.384581259397517
$ws.Range("G60").Value = 11.0468210526316
$ws.Range("H60").Value = 26.1187
$ws.Range("I60").Value = 24.27426
$ws.Range("L60").Value = 14.92165
$ws.Range("M60").Value = 21.4114
$ws.Range("N60").Value = 23.29473
$ws.Range("G61").Value = 11.0468210526316
$ws.Range("H61").Value = 26.1187
$ws.Range("I61").Value = 24.27426
$ws.Range("L61").Value = 14.92165
$ws.Range("M61").Value = 21.4114
$ws.Range("N61").Value = 23.29473
$ws.Range("G68").Value = 4441.67398852969
$ws.Range("H68").Value = 78333.4393117812
$ws.Range("G69").Value = 4441.67398852969
$ws.Range("H69").Value = 78333.4393117812
$ws.Range("G70").Value = 4441.67398852969
$ws.Range("H70").Value = 78333.4393117812
$ws.Range("G71").Value = 4441.67398852969
$ws.Range("H71").Value = 78333.4393117812
$ws.Range("G74").Value = 0.385883863094308
$ws.Range("G75").Value = 0.385883863094308
$ws.Range("F76").Value = 12.00235
$ws.Range("G76").Value = 11.608495
$ws.Range("I76").Value = 24.35405
$ws.Range("L76").Value = 15.01085
$ws.Range("M76").Value = 21.52653
$ws.Range("N76").Value = 23.95181
$ws.Range("F77").Value = 12.00235
$ws.Range("G77").Value = 11.608495
$ws.Range("I77").Value = 24.35405
$ws.Range("L77").Value = 15.01085
$ws.Range("M77").Value = 21.52653
$ws.Range("N77").Value = 23.95181
$ws.Range("G84").Value = 5641.44247552377
$ws.Range("H84").Value = 78333.4393117812
$ws.Range("I84").Value = 21947.55461
$ws.Range("N84").Value = 14000
$ws.Range("G85").Value = 5641.44247552377
$ws.Range("H85").Value = 78333.4393117812
$ws.Range("I85").Value = 21947.55461
$ws.Range("N85").Value = 14000
$ws.Range("G86").Value = 5641.44247552377
$ws.Range("H86").Value = 78333.4393117812
$ws.Range("I86").Value = 21947.55461
$ws.Range("N86").Value = 14000
$ws.Range("G87").Value = 5641.44247552377
$ws.Range("H87").Value = 78333.4393117812
$ws.Range("I87").Value = 21947.55461
$ws.Range("N87").Value = 14000
$ws.Range("G90").Value = 0.347623863094308
$ws.Range("G91").Value = 0.347623863094308
$ws.Range("G92").Value = 11.57813
$ws.Range("L92").Value = 15.52095
$ws.Range("M92").Value = 21.15237
$ws.Range("N92").Value = 23.91197
$ws.Range("G93").Value = 11.57813
$ws.Range("L93").Value = 15.52095
$ws.Range("M93").Value = 21.15237
$ws.Range("N93").Value = 23.91197
$ws.Range("G100").Value = 7785.36130343936
$ws.Range("H100").Value = 78333.4393117812
$ws.Range("I100").Value = 44045.56484
$ws.Range("M100").Value = 10230
$ws.Range("N100").Value = 33368.53277
$ws.Range("G101").Value = 7785.36130343936
$ws.Range("H101").Value = 78333.4393117812
$ws.Range("I101").Value = 44045.56484
$ws.Range("M101").Value = 10230
$ws.Range("N101").Value = 33368.53277
$ws.Range("G102").Value = 7785.36130343936
$ws.Range("H102").Value = 78333.4393117812
$ws.Range("I102").Value = 44045.56484
$ws.Range("M102").Value = 10230
$ws.Range("N102").Value = 33368.53277
$ws.Range("G103").Value = 7785.36130343936
$ws.Range("H103").Value = 78333.4393117812
$ws.Range("I103").Value = 44045.56484
$ws.Range("M103").Value = 10230
$ws.Range("N103").Value = 33368.53277
$ws.Range("G106").Value = 0.372889655108362
$ws.Range("G107").Value = 0.372889655108362
$ws.Range("G108").Value = 12.5766933333333
$ws.Range("G109").Value = 12.5766933333333
$ws.Range("G116").Value = 8224.18733695371
$ws.Range("H116").Value = 78333.4393117812
$ws.Range("I116").Value = 44437.91734
$ws.Range("M116").Value = 9862
$ws.Range("N116").Value = 36273.70655
$ws.Range("G117").Value = 8224.18733695371
$ws.Range("H117").Value = 78333.4393117812
$ws.Range("I117").Value = 44437.91734
$ws.Range("M117").Value = 9862
$ws.Range("N117").Value = 36273.70655
$ws.Range("G118").Value = 8224.18733695371
$ws.Range("H118").Value = 78333.4393117812
$ws.Range("I118").Value = 44437.91734
$ws.Range("M118").Value = 9862
$ws.Range("N118").Value = 36273.70655
$ws.Range("G119").Value = 8224.18733695371
$ws.Range("H119").Value = 78333.4393117812
$ws.Range("I119").Value = 44437.91734
$ws.Range("M119").Value = 9862
$ws.Range("N119").Value = 36273.70655
$ws.Range("G122").Value = 0.452465745979706
$ws.Range("G123").Value = 0.452465745979706
$ws.Range("G132").Value = 6006.74516471445
$ws.Range("H132").Value = 45353.4065086535
$ws.Range("I132").Value = 34579.02184
$ws.Range("N132").Value = 20061.11966
$ws.Range("G133").Value = 6006.74516471445
$ws.Range("H133").Value = 45353.4065086535
$ws.Range("I133").Value = 34579.02184
$ws.Range("N133").Value = 20061.11966
$ws.Range("G134").Value = 6006.74516471445
$ws.Range("H134").Value = 45353.4065086535
$ws.Range("I134").Value = 34579.02184
$ws.Range("N134").Value = 20061.11966
$ws.Range("G135").Value = 6006.74516471445
$ws.Range("H135").Value = 45353.4065086535
$ws.Range("I135").Value = 34579.02184
$ws.Range("N135").Value = 20061.11966
$ws.Range("G138").Value = 0.436250880015615
$ws.Range("G139").Value = 0.436250880015615
$ws.Range("G148").Value = 6319.91988678446
$ws.Range("H148").Value = 45353.4065086535
$ws.Range("I148").Value = 36394.75546
$ws.Range("N148").Value = 21475.94587
$ws.Range("G149").Value = 6319.91988678446
$ws.Range("H149").Value = 45353.4065086535
$ws.Range("I149").Value = 36394.75546
$ws.Range("N149").Value = 21475.94587
$ws.Range("G150").Value = 6319.91988678446
$ws.Range("H150").Value = 45353.4065086535
$ws.Range("I150").Value = 36394.75546
$ws.Range("N150").Value = 21475.94587
$ws.Range("G151").Value = 6319.91988678446
$ws.Range("H151").Value = 45353.4065086535
$ws.Range("I151").Value = 36394.75546
$ws.Range("N151").Value = 21475.94587
$ws.Range("G154").Value = 0.445795049428298
$ws.Range("G155").Value = 0.445795049428298

# --- Append new rows 162-177 (2019 - 2023 reporting period) ---
# Row 162: DRP (95th Percentile)
$ws.Range("A162").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B162").Value = "DRP (95th Percentile)"
$ws.Range("C162").Value = "D"
$ws.Range("D162").Value = "2019 - 2023"
$ws.Range("E162").Value = "Impact"
$ws.Range("F162").Value = 0.93
$ws.Range("G162").Value = 1.6854347826087
$ws.Range("H162").Value = 7.08
$ws.Range("I162").Value = 5.754
$ws.Range("J162").Value = $null
$ws.Range("K162").Value = $null
$ws.Range("L162").Value = 2.19
$ws.Range("M162").Value = 4.1056
$ws.Range("N162").Value = 4.9428
$ws.Range("O162").Value = 1805181.656
$ws.Range("P162").Value = 5545497.329
$ws.Range("Q162").Value = "Manawatu District"
$ws.Range("R162").Value = "Rangitīkei-Turakina"
$ws.Range("S162").Value = "Coastal Rangitikei"
$ws.Range("T162").Value = "Rang_4a"
$ws.Range("U162").Value = "mg/L"

# Row 163: DRP (Median)
$ws.Range("A163").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B163").Value = "DRP (Median)"
$ws.Range("C163").Value = "D"
$ws.Range("D163").Value = "2019 - 2023"
$ws.Range("E163").Value = "Impact"
$ws.Range("F163").Value = 0.93
$ws.Range("G163").Value = 1.6854347826087
$ws.Range("H163").Value = 7.08
$ws.Range("I163").Value = 5.754
$ws.Range("J163").Value = $null
$ws.Range("K163").Value = $null
$ws.Range("L163").Value = 2.19
$ws.Range("M163").Value = 4.1056
$ws.Range("N163").Value = 4.9428
$ws.Range("O163").Value = 1805181.656
$ws.Range("P163").Value = 5545497.329
$ws.Range("Q163").Value = "Manawatu District"
$ws.Range("R163").Value = "Rangitīkei-Turakina"
$ws.Range("S163").Value = "Coastal Rangitikei"
$ws.Range("T163").Value = "Rang_4a"
$ws.Range("U163").Value = "mg/L"

# Row 164: E coli (>260)
$ws.Range("A164").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B164").Value = "E coli (>260)"
$ws.Range("C164").Value = "E"
$ws.Range("D164").Value = "2019 - 2023"
$ws.Range("E164").Value = "Impact"
$ws.Range("F164").Value = 1593.5
$ws.Range("G164").Value = 5143.23489144265
$ws.Range("H164").Value = 45353.4065086535
$ws.Range("I164").Value = 38147.54463
$ws.Range("J164").Value = 73.9130434782609
$ws.Range("K164").Value = 86.9565217391304
$ws.Range("L164").Value = 1900
$ws.Range("M164").Value = 6040
$ws.Range("N164").Value = 10991.23377
$ws.Range("O164").Value = 1805181.656
$ws.Range("P164").Value = 5545497.329
$ws.Range("Q164").Value = "Manawatu District"
$ws.Range("R164").Value = "Rangitīkei-Turakina"
$ws.Range("S164").Value = "Coastal Rangitikei"
$ws.Range("T164").Value = "Rang_4a"
$ws.Range("U164").Value = "% exceedances over 260/100 mL"

# Row 165: E coli (>540)
$ws.Range("A165").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B165").Value = "E coli (>540)"
$ws.Range("C165").Value = "E"
$ws.Range("D165").Value = "2019 - 2023"
$ws.Range("E165").Value = "Impact"
$ws.Range("F165").Value = 1593.5
$ws.Range("G165").Value = 5143.23489144265
$ws.Range("H165").Value = 45353.4065086535
$ws.Range("I165").Value = 38147.54463
$ws.Range("J165").Value = 73.9130434782609
$ws.Range("K165").Value = 86.9565217391304
$ws.Range("L165").Value = 1900
$ws.Range("M165").Value = 6040
$ws.Range("N165").Value = 10991.23377
$ws.Range("O165").Value = 1805181.656
$ws.Range("P165").Value = 5545497.329
$ws.Range("Q165").Value = "Manawatu District"
$ws.Range("R165").Value = "Rangitīkei-Turakina"
$ws.Range("S165").Value = "Coastal Rangitikei"
$ws.Range("T165").Value = "Rang_4a"
$ws.Range("U165").Value = "% exceedances over 540/100 mL"

# Row 166: E coli (Median)
$ws.Range("A166").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B166").Value = "E coli (Median)"
$ws.Range("C166").Value = "E"
$ws.Range("D166").Value = "2019 - 2023"
$ws.Range("E166").Value = "Impact"
$ws.Range("F166").Value = 1593.5
$ws.Range("G166").Value = 5143.23489144265
$ws.Range("H166").Value = 45353.4065086535
$ws.Range("I166").Value = 38147.54463
$ws.Range("J166").Value = 73.9130434782609
$ws.Range("K166").Value = 86.9565217391304
$ws.Range("L166").Value = 1900
$ws.Range("M166").Value = 6040
$ws.Range("N166").Value = 10991.23377
$ws.Range("O166").Value = 1805181.656
$ws.Range("P166").Value = 5545497.329
$ws.Range("Q166").Value = "Manawatu District"
$ws.Range("R166").Value = "Rangitīkei-Turakina"
$ws.Range("S166").Value = "Coastal Rangitikei"
$ws.Range("T166").Value = "Rang_4a"
$ws.Range("U166").Value = "E. coli/100 mL"

# Row 167: E coli (95th Percentile)
$ws.Range("A167").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B167").Value = "E coli (95th Percentile)"
$ws.Range("C167").Value = "E"
$ws.Range("D167").Value = "2019 - 2023"
$ws.Range("E167").Value = "Impact"
$ws.Range("F167").Value = 1593.5
$ws.Range("G167").Value = 5143.23489144265
$ws.Range("H167").Value = 45353.4065086535
$ws.Range("I167").Value = 38147.54463
$ws.Range("J167").Value = 73.9130434782609
$ws.Range("K167").Value = 86.9565217391304
$ws.Range("L167").Value = 1900
$ws.Range("M167").Value = 6040
$ws.Range("N167").Value = 10991.23377
$ws.Range("O167").Value = 1805181.656
$ws.Range("P167").Value = 5545497.329
$ws.Range("Q167").Value = "Manawatu District"
$ws.Range("R167").Value = "Rangitīkei-Turakina"
$ws.Range("S167").Value = "Coastal Rangitikei"
$ws.Range("T167").Value = "Rang_4a"
$ws.Range("U167").Value = "E. coli/100 mL"

# Row 168: Ammoniacal-N (95th Percentile)
$ws.Range("A168").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B168").Value = "Ammoniacal-N (95th Percentile)"
$ws.Range("C168").Value = "D"
$ws.Range("D168").Value = "2019 - 2023"
$ws.Range("E168").Value = "Impact"
$ws.Range("F168").Value = 2.30782
$ws.Range("G168").Value = 5.06574087547498
$ws.Range("H168").Value = 29.7974519337602
$ws.Range("I168").Value = 23.53863
$ws.Range("J168").Value = $null
$ws.Range("K168").Value = $null
$ws.Range("L168").Value = 3.64094
$ws.Range("M168").Value = 10.42493
$ws.Range("N168").Value = 16.89208
$ws.Range("O168").Value = 1805181.656
$ws.Range("P168").Value = 5545497.329
$ws.Range("Q168").Value = "Manawatu District"
$ws.Range("R168").Value = "Rangitīkei-Turakina"
$ws.Range("S168").Value = "Coastal Rangitikei"
$ws.Range("T168").Value = "Rang_4a"
$ws.Range("U168").Value = "mg NH4-N/L"

# Row 169: Ammoniacal-N (Median)
$ws.Range("A169").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B169").Value = "Ammoniacal-N (Median)"
$ws.Range("C169").Value = "D"
$ws.Range("D169").Value = "2019 - 2023"
$ws.Range("E169").Value = "Impact"
$ws.Range("F169").Value = 2.30782
$ws.Range("G169").Value = 5.06574087547498
$ws.Range("H169").Value = 29.7974519337602
$ws.Range("I169").Value = 23.53863
$ws.Range("J169").Value = $null
$ws.Range("K169").Value = $null
$ws.Range("L169").Value = 3.64094
$ws.Range("M169").Value = 10.42493
$ws.Range("N169").Value = 16.89208
$ws.Range("O169").Value = 1805181.656
$ws.Range("P169").Value = 5545497.329
$ws.Range("Q169").Value = "Manawatu District"
$ws.Range("R169").Value = "Rangitīkei-Turakina"
$ws.Range("S169").Value = "Coastal Rangitikei"
$ws.Range("T169").Value = "Rang_4a"
$ws.Range("U169").Value = "mg NH4-N/L"

# Row 170: Nitrate-N (95th Percentile)
$ws.Range("A170").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B170").Value = "Nitrate-N (95th Percentile)"
$ws.Range("C170").Value = "B"
$ws.Range("D170").Value = "2019 - 2023"
$ws.Range("E170").Value = "Impact"
$ws.Range("F170").Value = 0.2185
$ws.Range("G170").Value = 0.491429691624456
$ws.Range("H170").Value = 3.89
$ws.Range("I170").Value = 1.87
$ws.Range("J170").Value = $null
$ws.Range("K170").Value = $null
$ws.Range("L170").Value = 0.122
$ws.Range("M170").Value = 0.74868
$ws.Range("N170").Value = 1.45
$ws.Range("O170").Value = 1805181.656
$ws.Range("P170").Value = 5545497.329
$ws.Range("Q170").Value = "Manawatu District"
$ws.Range("R170").Value = "Rangitīkei-Turakina"
$ws.Range("S170").Value = "Coastal Rangitikei"
$ws.Range("T170").Value = "Rang_4a"
$ws.Range("U170").Value = "mg NO3-N/L"

# Row 171: Nitrate-N (Median)
$ws.Range("A171").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B171").Value = "Nitrate-N (Median)"
$ws.Range("C171").Value = "A"
$ws.Range("D171").Value = "2019 - 2023"
$ws.Range("E171").Value = "Impact"
$ws.Range("F171").Value = 0.2185
$ws.Range("G171").Value = 0.491429691624456
$ws.Range("H171").Value = 3.89
$ws.Range("I171").Value = 1.87
$ws.Range("J171").Value = $null
$ws.Range("K171").Value = $null
$ws.Range("L171").Value = 0.122
$ws.Range("M171").Value = 0.74868
$ws.Range("N171").Value = 1.45
$ws.Range("O171").Value = 1805181.656
$ws.Range("P171").Value = 5545497.329
$ws.Range("Q171").Value = "Manawatu District"
$ws.Range("R171").Value = "Rangitīkei-Turakina"
$ws.Range("S171").Value = "Coastal Rangitikei"
$ws.Range("T171").Value = "Rang_4a"
$ws.Range("U171").Value = "mg NO3-N/L"

# Row 172: Soluble Inorganic Nitrogen (95th Percentile)
$ws.Range("A172").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B172").Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Range("C172").Value = $null
$ws.Range("D172").Value = "2019 - 2023"
$ws.Range("E172").Value = "Impact"
$ws.Range("F172").Value = 6.0285
$ws.Range("G172").Value = 9.14982608695652
$ws.Range("H172").Value = 38.74
$ws.Range("I172").Value = 29.54
$ws.Range("J172").Value = $null
$ws.Range("K172").Value = $null
$ws.Range("L172").Value = 10.42
$ws.Range("M172").Value = 13.8732
$ws.Range("N172").Value = 27.7974
$ws.Range("O172").Value = 1805181.656
$ws.Range("P172").Value = 5545497.329
$ws.Range("Q172").Value = "Manawatu District"
$ws.Range("R172").Value = "Rangitīkei-Turakina"
$ws.Range("S172").Value = "Coastal Rangitikei"
$ws.Range("T172").Value = "Rang_4a"
$ws.Range("U172").Value = "g/m3"

# Row 173: Soluble Inorganic Nitrogen (Median)
$ws.Range("A173").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B173").Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Range("C173").Value = $null
$ws.Range("D173").Value = "2019 - 2023"
$ws.Range("E173").Value = "Impact"
$ws.Range("F173").Value = 6.0285
$ws.Range("G173").Value = 9.14982608695652
$ws.Range("H173").Value = 38.74
$ws.Range("I173").Value = 29.54
$ws.Range("J173").Value = $null
$ws.Range("K173").Value = $null
$ws.Range("L173").Value = 10.42
$ws.Range("M173").Value = 13.8732
$ws.Range("N173").Value = 27.7974
$ws.Range("O173").Value = 1805181.656
$ws.Range("P173").Value = 5545497.329
$ws.Range("Q173").Value = "Manawatu District"
$ws.Range("R173").Value = "Rangitīkei-Turakina"
$ws.Range("S173").Value = "Coastal Rangitikei"
$ws.Range("T173").Value = "Rang_4a"
$ws.Range("U173").Value = "g/m3"

# Row 174: Total Nitrogen (95th Percentile)
$ws.Range("A174").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B174").Value = "Total Nitrogen (95th Percentile)"
$ws.Range("C174").Value = $null
$ws.Range("D174").Value = "2019 - 2023"
$ws.Range("E174").Value = "Impact"
$ws.Range("F174").Value = 8.37
$ws.Range("G174").Value = 13.1826086956522
$ws.Range("H174").Value = 47.1
$ws.Range("I174").Value = 40.2
$ws.Range("J174").Value = $null
$ws.Range("K174").Value = $null
$ws.Range("L174").Value = 16.4
$ws.Range("M174").Value = 20.304
$ws.Range("N174").Value = 38.594
$ws.Range("O174").Value = 1805181.656
$ws.Range("P174").Value = 5545497.329
$ws.Range("Q174").Value = "Manawatu District"
$ws.Range("R174").Value = "Rangitīkei-Turakina"
$ws.Range("S174").Value = "Coastal Rangitikei"
$ws.Range("T174").Value = "Rang_4a"
$ws.Range("U174").Value = "g/m3"

# Row 175: Total Nitrogen (Median)
$ws.Range("A175").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B175").Value = "Total Nitrogen (Median)"
$ws.Range("C175").Value = $null
$ws.Range("D175").Value = "2019 - 2023"
$ws.Range("E175").Value = "Impact"
$ws.Range("F175").Value = 8.37
$ws.Range("G175").Value = 13.1826086956522
$ws.Range("H175").Value = 47.1
$ws.Range("I175").Value = 40.2
$ws.Range("J175").Value = $null
$ws.Range("K175").Value = $null
$ws.Range("L175").Value = 16.4
$ws.Range("M175").Value = 20.304
$ws.Range("N175").Value = 38.594
$ws.Range("O175").Value = 1805181.656
$ws.Range("P175").Value = 5545497.329
$ws.Range("Q175").Value = "Manawatu District"
$ws.Range("R175").Value = "Rangitīkei-Turakina"
$ws.Range("S175").Value = "Coastal Rangitikei"
$ws.Range("T175").Value = "Rang_4a"
$ws.Range("U175").Value = "g/m3"

# Row 176: Total Phosphorus (95th Percentile)
$ws.Range("A176").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B176").Value = "Total Phosphorus (95th Percentile)"
$ws.Range("C176").Value = $null
$ws.Range("D176").Value = "2019 - 2023"
$ws.Range("E176").Value = "Impact"
$ws.Range("F176").Value = 1.765
$ws.Range("G176").Value = 2.80897826086957
$ws.Range("H176").Value = 9.66
$ws.Range("I176").Value = 9.126
$ws.Range("J176").Value = $null
$ws.Range("K176").Value = $null
$ws.Range("L176").Value = 3.485
$ws.Range("M176").Value = 5.4356
$ws.Range("N176").Value = 8.0258
$ws.Range("O176").Value = 1805181.656
$ws.Range("P176").Value = 5545497.329
$ws.Range("Q176").Value = "Manawatu District"
$ws.Range("R176").Value = "Rangitīkei-Turakina"
$ws.Range("S176").Value = "Coastal Rangitikei"
$ws.Range("T176").Value = "Rang_4a"
$ws.Range("U176").Value = "g/m3"

# Row 177: Total Phosphorus (Median)
$ws.Range("A177").Value = "Piakatutu at d/s Sanson STP"
$ws.Range("B177").Value = "Total Phosphorus (Median)"
$ws.Range("C177").Value = $null
$ws.Range("D177").Value = "2019 - 2023"
$ws.Range("E177").Value = "Impact"
$ws.Range("F177").Value = 1.765
$ws.Range("G177").Value = 2.80897826086957
$ws.Range("H177").Value = 9.66
$ws.Range("I177").Value = 9.126
$ws.Range("J177").Value = $null
$ws.Range("K177").Value = $null
$ws.Range("L177").Value = 3.485
$ws.Range("M177").Value = 5.4356
$ws.Range("N177").Value = 8.0258
$ws.Range("O177").Value = 1805181.656
$ws.Range("P177").Value = 5545497.329
$ws.Range("Q177").Value = "Manawatu District"
$ws.Range("R177").Value = "Rangitīkei-Turakina"
$ws.Range("S177").Value = "Coastal Rangitikei"
$ws.Range("T177").Value = "Rang_4a"
$ws.Range("U177").Value = "g/m3"

